$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 29, pushing the existing rows 29-48 down to 31-50.
$ws.Range("A29:A30").EntireRow.Insert()

# Row 29: new "Especial" quintal entry for the week of 2023-03-16 (serial 45001)
$ws.Cells.Item(29, 1).Value = 9
$ws.Cells.Item(29, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(29, 3).Value = 'Metropolitana'
$ws.Cells.Item(29, 4).Value = 45001
$ws.Cells.Item(29, 5).Value = 13
$ws.Cells.Item(29, 6).Value = 'Fruta'
$ws.Cells.Item(29, 7).Value = 100104
$ws.Cells.Item(29, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(29, 9).Value = 100104003
$ws.Cells.Item(29, 10).Value = 'Membrillo'
$ws.Cells.Item(29, 11).Value = 'Champion'
$ws.Cells.Item(29, 12).Value = 'Especial'
$ws.Cells.Item(29, 13).Value = 50
$ws.Cells.Item(29, 14).Value = 11000
$ws.Cells.Item(29, 15).Value = 11000
$ws.Cells.Item(29, 16).Value = 11000
$ws.Cells.Item(29, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(29, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(29, 19).Value = 611
$ws.Cells.Item(29, 20).Value = 18

# Row 30: new "Primera" quintal entry for the same week (2023-03-16, serial 45001)
$ws.Cells.Item(30, 1).Value = 9
$ws.Cells.Item(30, 2).Value = 'Vega Central Mapocho de Santiago'
$ws.Cells.Item(30, 3).Value = 'Metropolitana'
$ws.Cells.Item(30, 4).Value = 45001
$ws.Cells.Item(30, 5).Value = 13
$ws.Cells.Item(30, 6).Value = 'Fruta'
$ws.Cells.Item(30, 7).Value = 100104
$ws.Cells.Item(30, 8).Value = 'Frutos de pepita'
$ws.Cells.Item(30, 9).Value = 100104003
$ws.Cells.Item(30, 10).Value = 'Membrillo'
$ws.Cells.Item(30, 11).Value = 'Champion'
$ws.Cells.Item(30, 12).Value = 'Primera'
$ws.Cells.Item(30, 13).Value = 80
$ws.Cells.Item(30, 14).Value = 9000
$ws.Cells.Item(30, 15).Value = 9000
$ws.Cells.Item(30, 16).Value = 9000
$ws.Cells.Item(30, 17).Value = '$/caja 18 kilos granel'
$ws.Cells.Item(30, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(30, 19).Value = 500
$ws.Cells.Item(30, 20).Value = 18
